$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Cells.Item(2, 5).Value = '[''Normal'']'
$ws.Cells.Item(3, 4).Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Cells.Item(3, 5).Value = '[''Normal'', ''ParamViolation'']'
$ws.Cells.Item(11, 4).Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Cells.Item(11, 5).Value = '[''Normal'', ''HardwareFault'']'
$ws.Cells.Item(12, 4).Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Cells.Item(12, 5).Value = '[''Normal'', ''HardwareFault'']'
$ws.Cells.Item(15, 4).Value = '[0, 0, 0, 1, 0, 0, 0]'
$ws.Cells.Item(15, 5).Value = '[''ParamViolation'']'
$ws.Cells.Item(16, 4).Value = '[1, 0, 0, 0, 1, 0, 0]'
$ws.Cells.Item(16, 5).Value = '[''Normal'', ''RegulationViolation'']'
$ws.Cells.Item(19, 4).Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Cells.Item(19, 5).Value = '[]'
$ws.Cells.Item(24, 4).Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Cells.Item(24, 5).Value = '[''HardwareFault'']'
$ws.Cells.Item(25, 4).Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(25, 5).Value = '[''Normal'', ''SoftwareFault'']'
$ws.Cells.Item(26, 4).Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(26, 5).Value = '[''Normal'', ''SoftwareFault'']'
$ws.Cells.Item(27, 4).Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(27, 5).Value = '[''SoftwareFault'']'
$ws.Cells.Item(31, 4).Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(31, 5).Value = '[''Normal'', ''SoftwareFault'']'
$ws.Cells.Item(35, 4).Value = '[0, 0, 1, 0, 0, 0, 0]'
$ws.Cells.Item(35, 5).Value = '[''HardwareFault'']'
$ws.Cells.Item(38, 4).Value = '[1, 0, 1, 0, 0, 0, 0]'
$ws.Cells.Item(38, 5).Value = '[''Normal'', ''HardwareFault'']'
$ws.Cells.Item(54, 4).Value = '[0, 0, 0, 0, 0, 1, 0]'
$ws.Cells.Item(54, 5).Value = '[''CommunicationIssue'']'
$ws.Cells.Item(56, 4).Value = '[0, 0, 0, 0, 0, 0, 0]'
$ws.Cells.Item(56, 5).Value = '[]'
$ws.Cells.Item(58, 4).Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Cells.Item(58, 5).Value = '[''Normal'', ''ParamViolation'']'
$ws.Cells.Item(61, 4).Value = '[0, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(61, 5).Value = '[''SoftwareFault'']'
$ws.Cells.Item(68, 4).Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Cells.Item(68, 5).Value = '[''Normal'', ''ParamViolation'']'
$ws.Cells.Item(71, 4).Value = '[1, 0, 0, 1, 0, 0, 0]'
$ws.Cells.Item(71, 5).Value = '[''Normal'', ''ParamViolation'']'
$ws.Cells.Item(84, 4).Value = '[1, 0, 0, 0, 0, 0, 0]'
$ws.Cells.Item(84, 5).Value = '[''Normal'']'
$ws.Cells.Item(116, 4).Value = '[1, 0, 0, 0, 0, 0, 1]'
$ws.Cells.Item(116, 5).Value = '[''Normal'', ''SoftwareFault'']'
